# Update odds values in row 7 to reflect the latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = 1.44
$ws.Range("H7").Value = 4.33
$ws.Range("I7").Value = 8
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 9.5
$ws.Range("O7").Value = 1.29
$ws.Range("U7").Value = 2.2
$ws.Range("V7").Value = 1.62
$ws.Range("X7").Value = 6
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 9
$ws.Range("AB7").Value = 34
$ws.Range("AC7").Value = 9.5
$ws.Range("AH7").Value = 41
$ws.Range("AM7").Value = 900
$ws.Range("AO7").Value = 7
$ws.Range("AW7").Value = 8.5
$ws.Range("BA7").Value = 201
